$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the five course-day dates (2024 -> 2025 schedule)
$ws.Range("A2").Value  = "06/10/2025"
$ws.Range("A12").Value = "07/10/2025"
$ws.Range("A21").Value = "08/10/2025"
$ws.Range("A30").Value = "09/10/2025"
$ws.Range("A39").Value = "10/10/2025"

# Update the selected/active cell shown when the workbook is opened
$ws.Range("D16").Select()
